$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Thompson seedless" grapes at
# Terminal Hortofrutícola Agro Chillán; it belongs chronologically at the
# top of the existing Uva (K=Red Globe, ...) block, so insert a fresh row
# at row 73 and push the rest of the table (old rows 73-136) down to 74-137.
$ws.Rows.Item(73).Insert()

$ws.Cells.Item(73, 1).Value = 7
$ws.Cells.Item(73, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(73, 3).Value = "Ñuble"
$ws.Cells.Item(73, 4).Value = 44658
$ws.Cells.Item(73, 5).Value = 16
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100109
$ws.Cells.Item(73, 8).Value = "Uva"
$ws.Cells.Item(73, 9).Value = 100109001
$ws.Cells.Item(73, 10).Value = "Uva"
$ws.Cells.Item(73, 11).Value = "Thompson seedless"
$ws.Cells.Item(73, 12).Value = "Primera"
$ws.Cells.Item(73, 13).Value = 160
$ws.Cells.Item(73, 14).Value = 11000
$ws.Cells.Item(73, 15).Value = 12000
$ws.Cells.Item(73, 16).Value = 11500
$ws.Cells.Item(73, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(73, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(73, 19).Value = 639
$ws.Cells.Item(73, 20).Value = 18
